$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.588.40'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '3.361.07'
$ws.Range("E3").Value = '  -2.07%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.00'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.58'
$ws.Range("E6").Value = '  +0.86%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("D8").Value = '3.352.60'
$ws.Range("E8").Value = '  -2.04%  '

$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.629'
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.163'
$ws.Range("E11").Value = '  +3.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.45'
$ws.Range("E12").Value = '  -1.37%  '

$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("E14").Value = '  -0.37%  '

$ws.Range("D15").Value = '3.897.25'
$ws.Range("E15").Value = '  -2.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.42'
$ws.Range("E16").Value = '  +1.97%  '

$ws.Range("E17").Value = '  -2.08%  '

$ws.Range("D18").Value = '3.367.79'
$ws.Range("E18").Value = '  -1.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.84'
$ws.Range("E19").Value = '  -0.12%  '

$ws.Range("D20").Value = '64.519.21'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.986'
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '464.11'
$ws.Range("E22").Value = '  +13.78%  '

$ws.Range("E23").Value = '  +10.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.08'
$ws.Range("E24").Value = '  -2.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.09'
$ws.Range("E25").Value = '  +3.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.32'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").Value = '  +0.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.84'
$ws.Range("E28").Value = '  +1.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.78'
$ws.Range("E29").Value = '  -1.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.14'
$ws.Range("E30").Value = '  +1.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.64'
$ws.Range("E31").Value = '  -0.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.49'
$ws.Range("E32").Value = '  -0.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '582.99'
$ws.Range("E33").Value = '  -0.50%  '

$ws.Range("E34").Value = '  +0.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.94'
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("E36").Value = '  +0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.139'
$ws.Range("E37").Value = '  -8.80%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.49'
$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("E39").Value = '  -0.59%  '

$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.375'
$ws.Range("E41").Value = '  -0.33%  '

$ws.Range("D42").Value = '3.104.18'
$ws.Range("E42").Value = '  -2.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'
$ws.Range("E44").Value = '  -3.84%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.53'
$ws.Range("E45").Value = '  +1.36%  '

$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.22'
$ws.Range("E47").Value = '  -0.13%  '

$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("E49").Value = '  -2.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.38'
$ws.Range("E50").Value = '  -0.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.37'
$ws.Range("E51").Value = '  -0.76%  '
